# This script refreshes the Universalis market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ -- columns H:N)
# on each job sheet, matching the latest scheduled market-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 704.3043
$ws.Range("I19").Value = 645.2222
$ws.Range("J19").Value = 742.2857
$ws.Range("K19").Value = 645.2222
$ws.Range("L19").Value = 742.2857
$ws.Range("M19").Value = -470.2222
$ws.Range("N19").Value = -1092.2857
$ws.Range("H39").Value = 140.40909
$ws.Range("I39").Value = 67.083336
$ws.Range("J39").Value = 228.4
$ws.Range("K39").Value = 201.250008
$ws.Range("L39").Value = 685.2
$ws.Range("M39").Value = 94.74999199999999
$ws.Range("N39").Value = -1277.2
$ws.Range("H40").Value = 3237
$ws.Range("J40").Value = 2461.2
$ws.Range("L40").Value = 2461.2
$ws.Range("N40").Value = -2811.2
$ws.Range("H133").Value = 14766.667
$ws.Range("J133").Value = 14766.667
$ws.Range("L133").Value = 14766.667
$ws.Range("N133").Value = -24886.667
$ws.Range("H135").Value = 6784.95
$ws.Range("I135").Value = 7927.75
$ws.Range("J135").Value = 2213.75
$ws.Range("K135").Value = 71349.75
$ws.Range("L135").Value = 19923.75
$ws.Range("M135").Value = -68814.75
$ws.Range("N135").Value = -24993.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = -2885
$ws.Range("H102").Value = 2533.3333
$ws.Range("I102").Value = 2533.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2533.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -911.3332999999998
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2111.1765
$ws.Range("I122").Value = 1944.5714
$ws.Range("K122").Value = 5833.7142
$ws.Range("M122").Value = -3383.7142
$ws.Range("H133").Value = 53000
$ws.Range("J133").Value = 53000
$ws.Range("L133").Value = 53000
$ws.Range("N133").Value = -58060
$ws.Range("H139").Value = 65333.332
$ws.Range("J139").Value = 65333.332
$ws.Range("L139").Value = 65333.332
$ws.Range("N139").Value = -75613.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 49800
$ws.Range("J59").Value = 49800
$ws.Range("L59").Value = 49800
$ws.Range("N59").Value = -51494
$ws.Range("H86").Value = 11156.125
$ws.Range("I86").Value = 5730.077
$ws.Range("J86").Value = 34669
$ws.Range("K86").Value = 5730.077
$ws.Range("L86").Value = 34669
$ws.Range("M86").Value = -4607.077
$ws.Range("N86").Value = -36915
$ws.Range("H89").Value = 11156.125
$ws.Range("I89").Value = 5730.077
$ws.Range("J89").Value = 34669
$ws.Range("K89").Value = 28650.385
$ws.Range("L89").Value = 173345
$ws.Range("M89").Value = -23034.385
$ws.Range("N89").Value = -184577
$ws.Range("H99").Value = 2478.0908
$ws.Range("I99").Value = 709.875
$ws.Range("J99").Value = 7193.3335
$ws.Range("K99").Value = 709.875
$ws.Range("L99").Value = 7193.3335
$ws.Range("M99").Value = 788.125
$ws.Range("N99").Value = -10189.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1258.1538
$ws.Range("I58").Value = 652.8077
$ws.Range("J58").Value = 2468.8462
$ws.Range("K58").Value = 652.8077
$ws.Range("L58").Value = 2468.8462
$ws.Range("M58").Value = -449.8077
$ws.Range("N58").Value = -2874.8462
$ws.Range("H86").Value = 1497.5714
$ws.Range("I86").Value = 1320.8889
$ws.Range("J86").Value = 1815.6
$ws.Range("K86").Value = 1320.8889
$ws.Range("L86").Value = 1815.6
$ws.Range("M86").Value = -197.8888999999999
$ws.Range("N86").Value = -4061.6
$ws.Range("H89").Value = 1497.5714
$ws.Range("I89").Value = 1320.8889
$ws.Range("J89").Value = 1815.6
$ws.Range("K89").Value = 6604.4445
$ws.Range("L89").Value = 9078
$ws.Range("M89").Value = -988.4444999999996
$ws.Range("N89").Value = -20310
$ws.Range("H103").Value = 15857.833
$ws.Range("I103").Value = 12069.4
$ws.Range("K103").Value = 12069.4
$ws.Range("M103").Value = -10897.4
$ws.Range("H134").Value = 5817.6665
$ws.Range("J134").Value = 8083.5
$ws.Range("L134").Value = 24250.5
$ws.Range("N134").Value = -29320.5
$ws.Range("H136").Value = 1258.1538
$ws.Range("I136").Value = 652.8077
$ws.Range("J136").Value = 2468.8462
$ws.Range("K136").Value = 1958.4231
$ws.Range("L136").Value = 7406.5386
$ws.Range("M136").Value = 591.5769
$ws.Range("N136").Value = -12506.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 542
$ws.Range("I36").Value = 491.42856
$ws.Range("J36").Value = 660
$ws.Range("K36").Value = 1474.28568
$ws.Range("L36").Value = 1980
$ws.Range("M36").Value = -1305.28568
$ws.Range("N36").Value = -2318
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -314
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 1428
$ws.Range("N89").Value = -16356
$ws.Range("H106").Value = 3573.375
$ws.Range("J106").Value = 3726.7144
$ws.Range("L106").Value = 11180.1432
$ws.Range("N106").Value = -13072.1432
$ws.Range("H107").Value = 245190.1
$ws.Range("J107").Value = 386413.62
$ws.Range("L107").Value = 1159240.86
$ws.Range("N107").Value = -1163080.86
$ws.Range("H115").Value = 2531.6667
$ws.Range("I115").Value = 1200
$ws.Range("J115").Value = 2798
$ws.Range("K115").Value = 3600
$ws.Range("L115").Value = 8394
$ws.Range("M115").Value = -2425
$ws.Range("N115").Value = -10744
$ws.Range("H129").Value = 1382.9
$ws.Range("I129").Value = 315
$ws.Range("J129").Value = 1649.875
$ws.Range("K129").Value = 945
$ws.Range("L129").Value = 4949.625
$ws.Range("M129").Value = 4055
$ws.Range("N129").Value = -14949.625
$ws.Range("H131").Value = 2292.5618
$ws.Range("J131").Value = 2534.9114
$ws.Range("L131").Value = 7604.7342
$ws.Range("N131").Value = -17684.7342

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1591587.2
$ws.Range("I122").Value = 2227022.2
$ws.Range("K122").Value = 6681066.600000001
$ws.Range("M122").Value = -6678616.600000001
$ws.Range("H138").Value = 62500
$ws.Range("J138").Value = 62500
$ws.Range("L138").Value = 62500
$ws.Range("N138").Value = -72780
$ws.Range("H139").Value = 59779.4
$ws.Range("J139").Value = 59779.4
$ws.Range("L139").Value = 59779.4
$ws.Range("N139").Value = -70059.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3517.8572
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 3480.7693
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 3480.7693
$ws.Range("M40").Value = -3864
$ws.Range("N40").Value = -3752.7693
$ws.Range("H55").Value = 228.125
$ws.Range("J55").Value = 381
$ws.Range("L55").Value = 381
$ws.Range("N55").Value = -727
$ws.Range("H68").Value = 1768.1428
$ws.Range("I68").Value = 1585.5714
$ws.Range("J68").Value = 1950.7142
$ws.Range("K68").Value = 1585.5714
$ws.Range("L68").Value = 1950.7142
$ws.Range("M68").Value = -836.5714
$ws.Range("N68").Value = -3448.7142
$ws.Range("H71").Value = 1768.1428
$ws.Range("I71").Value = 1585.5714
$ws.Range("J71").Value = 1950.7142
$ws.Range("K71").Value = 7927.857
$ws.Range("L71").Value = 9753.571
$ws.Range("M71").Value = -4183.857
$ws.Range("N71").Value = -17241.571
$ws.Range("H122").Value = 3644.111
$ws.Range("I122").Value = 2398.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7195.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -4745.5
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2629.8333
$ws.Range("I96").Value = 2444.75
$ws.Range("K96").Value = 2444.75
$ws.Range("M96").Value = -1071.75
